$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "line_pair"
$ws.Range("F1").Value = "direction"

$ws.Range("D2:D7").Value = 100

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = -1

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1

$ws.Range("E5").Value = 2
$ws.Range("F5").Value = -1

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = -1

$ws.Range("H4").Select()
